# Apply new values to the B2:F8 data block on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @{
    "B2" = -0.1562
    "C2" = -0.1562
    "D2" = -0.1643
    "E2" = 0.1111
    "F2" = -0.5402

    "B3" = -0.1712
    "C3" = -0.1712
    "D3" = -0.1764
    "E3" = 0.2046
    "F3" = -0.4198

    "B4" = -0.0896
    "C4" = -0.0896
    "D4" = -0.08069999999999999
    "E4" = 0.1335
    "F4" = -0.3636

    "B5" = 0.0517
    "C5" = 0.0517
    "D5" = 0.0586
    "E5" = 0.2607
    "F5" = -0.2958

    "B6" = 0.0512
    "C6" = 0.0512
    "D6" = 0.0456
    "E6" = 0.2034
    "F6" = -0.3402

    "B7" = 0.0503
    "C7" = 0.0503
    "D7" = 0.0444
    "E7" = 0.1801
    "F7" = -0.1017

    "B8" = 0.0842
    "C8" = 0.0842
    "D8" = 0.07829999999999999
    "E8" = 0.214
    "F8" = -0.143
}

foreach ($addr in $data.Keys) {
    $ws.Range($addr).Value = $data[$addr]
}
